$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1111.7778
$ws.Range("I98").Value = 959
$ws.Range("J98").Value = 1302.75
$ws.Range("K98").Value = 959
$ws.Range("L98").Value = 1302.75
$ws.Range("M98").Value = 539
$ws.Range("N98").Value = -4298.75

$ws.Range("H122").Value = 1111.7778
$ws.Range("I122").Value = 959
$ws.Range("J122").Value = 1302.75
$ws.Range("K122").Value = 2877
$ws.Range("L122").Value = 3908.25
$ws.Range("M122").Value = -427
$ws.Range("N122").Value = -8808.25

$ws.Range("H127").Value = 1250
$ws.Range("J127").Value = 1250
$ws.Range("L127").Value = 3750
$ws.Range("N127").Value = -13670

$ws.Range("H129").Value = 1929.75
$ws.Range("I129").Value = 1219.6
$ws.Range("K129").Value = 3658.8
$ws.Range("M129").Value = 1341.2

$ws.Range("H134").Value = 275000
$ws.Range("J134").Value = 275000
$ws.Range("L134").Value = 275000
$ws.Range("N134").Value = -285140

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3991.8
$ws.Range("I2").Value = 3766.5
$ws.Range("J2").Value = 4442.4
$ws.Range("K2").Value = 3766.5
$ws.Range("L2").Value = 4442.4
$ws.Range("M2").Value = -3653.5
$ws.Range("N2").Value = -4668.4

$ws.Range("H52").Value = 83500
$ws.Range("J52").Value = 83500
$ws.Range("L52").Value = 83500
$ws.Range("N52").Value = -84136

$ws.Range("H74").Value = 7680.2
$ws.Range("J74").Value = 8477
$ws.Range("L74").Value = 8477
$ws.Range("N74").Value = -10225

$ws.Range("H77").Value = 7680.2
$ws.Range("J77").Value = 8477
$ws.Range("L77").Value = 42385
$ws.Range("N77").Value = -51121

$ws.Range("H88").Value = 997.5
$ws.Range("I88").Value = 455
$ws.Range("J88").Value = 1178.3334
$ws.Range("K88").Value = 455
$ws.Range("L88").Value = 1178.3334
$ws.Range("M88").Value = -49
$ws.Range("N88").Value = -1990.3334

$ws.Range("H91").Value = 997.5
$ws.Range("I91").Value = 455
$ws.Range("J91").Value = 1178.3334
$ws.Range("K91").Value = 455
$ws.Range("L91").Value = 1178.3334
$ws.Range("M91").Value = 949
$ws.Range("N91").Value = -3986.3334

$ws.Range("H97").Value = 989.3158
$ws.Range("I97").Value = 737.53845
$ws.Range("K97").Value = 737.53845
$ws.Range("M97").Value = -241.53845

$ws.Range("H102").Value = 3558.3076
$ws.Range("I102").Value = 695.3333
$ws.Range("K102").Value = 695.3333
$ws.Range("M102").Value = 926.6667

$ws.Range("H104").Value = 30316.666
$ws.Range("J104").Value = 30316.666
$ws.Range("L104").Value = 30316.666
$ws.Range("N104").Value = -37304.666

$ws.Range("H116").Value = 3991.8
$ws.Range("I116").Value = 3766.5
$ws.Range("J116").Value = 4442.4
$ws.Range("K116").Value = 3766.5
$ws.Range("L116").Value = 4442.4
$ws.Range("M116").Value = -1472.5
$ws.Range("N116").Value = -9030.4

$ws.Range("H122").Value = 2251.5715
$ws.Range("I122").Value = 1152.2
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 3456.6
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -1006.6
$ws.Range("N122").Value = -19900

$ws.Range("H132").Value = 1000
$ws.Range("I132").Value = 1000
$ws.Range("K132").Value = 3000
$ws.Range("M132").Value = -470

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3991.8
$ws.Range("I3").Value = 3766.5
$ws.Range("J3").Value = 4442.4
$ws.Range("K3").Value = 3766.5
$ws.Range("L3").Value = 4442.4
$ws.Range("M3").Value = -3652.5
$ws.Range("N3").Value = -4670.4

$ws.Range("H105").Value = 1206.4286
$ws.Range("I105").Value = 1149
$ws.Range("K105").Value = 1149
$ws.Range("M105").Value = 598

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2500.1428
$ws.Range("I99").Value = 1800.4
$ws.Range("K99").Value = 1800.4
$ws.Range("M99").Value = -302.4000000000001

$ws.Range("H122").Value = 1351.1666
$ws.Range("I122").Value = 1302
$ws.Range("J122").Value = 1449.5
$ws.Range("K122").Value = 3906
$ws.Range("L122").Value = 4348.5
$ws.Range("M122").Value = -1456
$ws.Range("N122").Value = -9248.5

$ws.Range("H126").Value = 2500.1428
$ws.Range("I126").Value = 1800.4
$ws.Range("K126").Value = 5401.200000000001
$ws.Range("M126").Value = -2931.200000000001

$ws.Range("H141").Value = 80625.3
$ws.Range("J141").Value = 80625.3
$ws.Range("L141").Value = 80625.3
$ws.Range("N141").Value = -90985.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 5000
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 5000
$ws.Range("K31").Value = 0
$ws.Range("L31").ClearContents()
$ws.Range("M31").Value = 15000
$ws.Range("N31").Value = -15576

$ws.Range("H122").Value = 764.5
$ws.Range("I122").Value = 745.5
$ws.Range("J122").Value = 802.5
$ws.Range("K122").Value = 6709.5
$ws.Range("L122").Value = 7222.5
$ws.Range("M122").Value = -4259.5
$ws.Range("N122").Value = -12122.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 6250
$ws.Range("I40").Value = 5000
$ws.Range("K40").Value = 5000
$ws.Range("M40").Value = -4849

$ws.Range("H43").Value = 3496.6667
$ws.Range("I43").Value = 3496.6667
$ws.Range("K43").Value = 3496.6667
$ws.Range("M43").Value = -3345.6667

$ws.Range("H70").Value = 4101.8
$ws.Range("I70").Value = 2500
$ws.Range("K70").Value = 2500
$ws.Range("M70").Value = -2230

$ws.Range("H73").Value = 4101.8
$ws.Range("I73").Value = 2500
$ws.Range("K73").Value = 2500
$ws.Range("M73").Value = -1564

$ws.Range("H97").Value = 964.875
$ws.Range("I97").Value = 541
$ws.Range("J97").Value = 1219.2
$ws.Range("K97").Value = 541
$ws.Range("L97").Value = 1219.2
$ws.Range("M97").Value = -45
$ws.Range("N97").Value = -2211.2

$ws.Range("H122").Value = 2679.7368
$ws.Range("I122").Value = 2054.5293
$ws.Range("K122").Value = 6163.5879
$ws.Range("M122").Value = -3713.5879

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

$ws.Range("H7").Value = 8331
$ws.Range("I7").Value = 7500
$ws.Range("J7").Value = 9162
$ws.Range("K7").Value = 7500
$ws.Range("L7").Value = 9162
$ws.Range("M7").Value = -7388
$ws.Range("N7").Value = -9386

$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("M28").ClearContents()

$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()

$ws.Range("H40").Value = 3649.6667
$ws.Range("I40").Value = 3649.6667
$ws.Range("K40").Value = 3649.6667
$ws.Range("M40").Value = -3513.6667

$ws.Range("H126").Value = 8331
$ws.Range("I126").Value = 7500
$ws.Range("J126").Value = 9162
$ws.Range("K126").Value = 22500
$ws.Range("L126").Value = 27486
$ws.Range("M126").Value = -20030
$ws.Range("N126").Value = -32426

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 524.7143
$ws.Range("I81").Value = 537.1667
$ws.Range("J81").Value = 450
$ws.Range("K81").Value = 1074.3334
$ws.Range("L81").Value = 900
$ws.Range("M81").Value = -13.33339999999998
$ws.Range("N81").Value = -3022

$ws.Range("H84").Value = 524.7143
$ws.Range("I84").Value = 537.1667
$ws.Range("J84").Value = 450
$ws.Range("K84").Value = 5371.666999999999
$ws.Range("L84").Value = 4500
$ws.Range("M84").Value = -67.66699999999946
$ws.Range("N84").Value = -15108

$ws.Range("H122").Value = 973.7826
$ws.Range("I122").Value = 954.85
$ws.Range("J122").Value = 1100
$ws.Range("K122").Value = 2864.55
$ws.Range("L122").Value = 3300
$ws.Range("M122").Value = -414.5500000000002
$ws.Range("N122").Value = -8200

$ws.Range("H126").Value = 5834.625
$ws.Range("J126").Value = 7483.3335
$ws.Range("L126").Value = 22450.0005
$ws.Range("N126").Value = -27390.0005
